$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-15 changes from 45243 to 45244 (one day later)
$ws.Range("C2:C15").Value = 45244
